{"js": "// Locate the paragraph with the exact text \"Farmer's Stories\" (curly apostrophe)\n// and set its font color to the green used elsewhere in this doc for\n// \"handled\" index items (00864B), matching both the run and the paragraph\n// mark formatting.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst target = \"Farmer\\u2019s Stories\";\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  para.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text === target) {\n    para.font.color = \"#00864B\";\n    break;\n  }\n}\nawait context.sync();\n", "ps1": "# Find the paragraph whose text is \"Farmer's Stories\" (curly apostrophe)\n# and color both its run text and paragraph mark with the green (00864B)\n# used elsewhere in the index for handled items.\n$d = $word.ActiveDocument\n\n$target = [string][char]0x2019\n$target = \"Farmer\" + $target + \"s Stories\"\n\n# Word COM colors are packed as 0x00BBGGRR (the classic VB RGB() order),\n# so build the value from the RGB triple 0x00, 0x86, 0x4B rather than\n# using the OOXML RRGGBB hex directly.\n$wdColor = 0x00 + (0x86 * 0x100) + (0x4B * 0x10000)\n\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    $text = $r.Text.TrimEnd([char]13, [char]7)\n    if ($text -eq $target) {\n        $r.Font.Color = $wdColor\n        break\n    }\n}\n"}
